$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numara / Ad Soyad / Bölüm block (I7:K9 merged cells already exist)
$ws.Range("I7").Value = 20215070019
$ws.Range("I8").Value = "KÜBRA ÇABUK"
$ws.Range("I9").Value = "YBS"

# Lookup demo block
$ws.Range("C18").Value = 333
$ws.Range("D18").Formula = '=VLOOKUP($C$18,$C$3:$F$15,2)'
$ws.Range("E18").Formula = '=VLOOKUP($C$18,$C$3:$F$15,3)'
$ws.Range("F18").Formula = '=VLOOKUP($C$18,$C$3:$F$15,4)'
# F18 picks up the same (thinner) right border as D18/E18 instead of the
# heavier outer-table edge it had before
$ws.Range("F18").Borders.Item(10).Weight = 2

# New cell below with a comma
$ws.Range("D19").Value = ","

# Move active selection
$ws.Range("F19").Select()
